$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 ("Kund") renumbering of "Main flow" steps ---
# Replace the old "N.1 / N.2 / ..." prefixed numbering with a simple "1 / 2 / 3 / 4" numbering.
$ws1.Range("B14").Value = "1 Kunden väljer att lägga till produkt från menyn`n2 Produkten läggs till i varukorgen"
$ws1.Range("B20").Value = "1 Kunden väljer att lägga till produkt från menyn`n2 Produkten läggs till i varukorgen"
$ws1.Range("B26").Value = "1 Kunden väljer att minska antal`n2 Antalet produkter minskas`n3 Varukorgen uppdateras`n4 Varukorgen visas på nytt"
$ws1.Range("B33").Value = "1 Kunden väljer att ta bort varan`n2 Varan tas bort från varukorgen`n3 Varukorgen uppdateras`n4 Varukorgen visas på nytt"
$ws1.Range("B39").Value = "1 Kunden väljer att ändra produkt`n2 Produkten tas bort från varukorgen`n3 Kunden tas till menyn för välj produkt"
$ws1.Range("B45").Value = "1 Kunden väljer att öka antal`n2 Antalet produkter ökas`n3 Varukorgen uppdateras`n4 Varukorgen visas på nytt"
$ws1.Range("B51").Value = "1 Kunden väljer att minska antal`n2 Antalet produkter minskas`n3 Varukorgen uppdateras`n4 Varukorgen visas på nytt"

# --- Sheet2 ("Ägare") table: pre-conditions & main flow rewritten (owner table) ---
$ws2.Range("B5").Value = "* app för ägare finns installerad`n* produkten måste finnas skapad"
$ws2.Range("B7").Value = "1. Ägaren väljer ""lägg till produkt"" `n2. Lista med tillgängliga produkter visas`n3. produkt väljs`n4. produkt läggs till på menyn"

# Row heights for the wrapped cells shrink because the text now has fewer lines.
$ws2.Rows.Item(5).RowHeight = 28.8
$ws2.Rows.Item(7).RowHeight = 58.2

# --- View state: update selection/scroll on sheet1, then restore sheet2 as the active tab ---
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("B52").Select()

$ws2.Activate()
$ws2.Range("F14").Select()
